function RGB($r, $g, $b) {
  return $r + ($g * 256) + ($b * 65536)
}

# The deck's single reachable theme (ppt/theme/theme1.xml, used by the
# slide master / all slides) currently carries the "Integral" color
# scheme. The edit swaps it for the stock "Office Theme" color scheme
# (font scheme / format scheme are already identical between the two
# themes embedded in this deck, so only the 12 theme colors change).

$p   = $ppt.ActivePresentation
$s   = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$tcs.Colors(1).RGB  = RGB 0x00 0x00 0x00   # Dark 1   (dk1)
$tcs.Colors(2).RGB  = RGB 0xFF 0xFF 0xFF   # Light 1  (lt1)
$tcs.Colors(3).RGB  = RGB 0x44 0x54 0x6A   # Dark 2   (dk2)
$tcs.Colors(4).RGB  = RGB 0xE7 0xE6 0xE6   # Light 2  (lt2)
$tcs.Colors(5).RGB  = RGB 0x5B 0x9B 0xD5   # Accent 1
$tcs.Colors(6).RGB  = RGB 0xED 0x7D 0x31   # Accent 2
$tcs.Colors(7).RGB  = RGB 0xA5 0xA5 0xA5   # Accent 3
$tcs.Colors(8).RGB  = RGB 0xFF 0xC0 0x00   # Accent 4
$tcs.Colors(9).RGB  = RGB 0x44 0x72 0xC4   # Accent 5
$tcs.Colors(10).RGB = RGB 0x70 0xAD 0x47   # Accent 6
$tcs.Colors(11).RGB = RGB 0x05 0x63 0xC1   # Hyperlink
$tcs.Colors(12).RGB = RGB 0x95 0x4F 0x72   # Followed Hyperlink
